$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt9a"
$ws.Range("C2").Value = "Fzd10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.103667
$ws.Range("H2").Value = 0.311001
$ws.Range("I2").Value = 0.01587706838992035
$ws.Range("J2").Value = 0.01587706838992035
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.06694666666666667
$ws.Range("N2").Value = 0.20084
$ws.Range("O2").Value = 0.5098924310779488
$ws.Range("P2").Value = 0.5098924310779488
$ws.Range("Q2").Value = 0.006940160093333334
$ws.Range("R2").Value = 0.06246144084000001
$ws.Range("S2").Value = 0.008095596999727344
$ws.Range("T2").Value = 0.008095596999727343

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt9a"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.103667
$ws.Range("H3").Value = 0.311001
$ws.Range("I3").Value = 0.01587706838992035
$ws.Range("J3").Value = 0.01587706838992035
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.064349
$ws.Range("N3").Value = 0.193047
$ws.Range("O3").Value = 0.4901075689220513
$ws.Range("P3").Value = 0.4901075689220513
$ws.Range("Q3").Value = 0.006670867783000001
$ws.Range("R3").Value = 0.06003781004700001
$ws.Range("S3").Value = 0.007781471390193011
$ws.Range("T3").Value = 0.00778147139019301

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt9a"
$ws.Range("C4").Value = "Fzd10"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.232908333333333
$ws.Range("H4").Value = 9.698725
$ws.Range("I4").Value = 0.4951344854840667
$ws.Range("J4").Value = 0.4951344854840667
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.06694666666666667
$ws.Range("N4").Value = 0.20084
$ws.Range("O4").Value = 0.5098924310779488
$ws.Range("P4").Value = 0.5098924310779488
$ws.Range("Q4").Value = 0.2164324365555556
$ws.Range("R4").Value = 1.947891929
$ws.Range("S4").Value = 0.2524653265140002
$ws.Range("T4").Value = 0.2524653265140002

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt9a"
$ws.Range("C5").Value = "Fzd10"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.232908333333333
$ws.Range("H5").Value = 9.698725
$ws.Range("I5").Value = 0.4951344854840667
$ws.Range("J5").Value = 0.4951344854840667
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.064349
$ws.Range("N5").Value = 0.193047
$ws.Range("O5").Value = 0.4901075689220513
$ws.Range("P5").Value = 0.4901075689220513
$ws.Range("Q5").Value = 0.2080344183416667
$ws.Range("R5").Value = 1.872309765075
$ws.Range("S5").Value = 0.2426691589700666
$ws.Range("T5").Value = 0.2426691589700666

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Wnt9a"
$ws.Range("C6").Value = "Fzd10"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.192778666666667
$ws.Range("H6").Value = 9.578336
$ws.Range("I6").Value = 0.4889884461260129
$ws.Range("J6").Value = 0.4889884461260129
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.06694666666666667
$ws.Range("N6").Value = 0.20084
$ws.Range("O6").Value = 0.5098924310779488
$ws.Range("P6").Value = 0.5098924310779488
$ws.Range("Q6").Value = 0.2137458891377778
$ws.Range("R6").Value = 1.92371300224
$ws.Range("S6").Value = 0.2493315075642213
$ws.Range("T6").Value = 0.2493315075642213

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Wnt9a"
$ws.Range("C7").Value = "Fzd10"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.192778666666667
$ws.Range("H7").Value = 9.578336
$ws.Range("I7").Value = 0.4889884461260129
$ws.Range("J7").Value = 0.4889884461260129
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.064349
$ws.Range("N7").Value = 0.193047
$ws.Range("O7").Value = 0.4901075689220513
$ws.Range("P7").Value = 0.4901075689220513
$ws.Range("Q7").Value = 0.2054521144213333
$ws.Range("R7").Value = 1.849069029792
$ws.Range("S7").Value = 0.2396569385617916
$ws.Range("T7").Value = 0.2396569385617916
